$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D8","D9","D10","D11","D14","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D30","D31","D34","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D50")
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.261.31'
$ws.Range("E2").Value = '  +1.58%  '

$ws.Range("D3").Value = '2.646.74'
$ws.Range("E3").Value = '  +1.93%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '582.19'
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("D6").Value = '144.29'
$ws.Range("E6").Value = '  +1.02%  '

$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("D9").Value = '6.56'
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").Value = '0.109'
$ws.Range("E10").Value = '  +2.84%  '

$ws.Range("D11").Value = '0.159'
$ws.Range("E11").Value = '  +1.54%  '

$ws.Range("E12").Value = '  +2.47%  '

$ws.Range("D13").Value = '3.112.94'
$ws.Range("E13").Value = '  +1.92%  '

$ws.Range("D14").Value = '26.31'
$ws.Range("E14").Value = '  +6.44%  '

$ws.Range("D15").Value = '61.208.44'
$ws.Range("E15").Value = '  +1.48%  '

$ws.Range("D16").Value = '0.0000146'
$ws.Range("E16").Value = '  +2.61%  '

$ws.Range("D17").Value = '2.655.67'
$ws.Range("E17").Value = '  +2.02%  '

$ws.Range("D18").Value = '11.70'
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("D19").Value = '4.76'
$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("D20").Value = '353.23'
$ws.Range("E20").Value = '  +1.58%  '

$ws.Range("D21").Value = '6.88'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").Value = '0.527'
$ws.Range("E23").Value = '  +0.68%  '

$ws.Range("D24").Value = '64.38'
$ws.Range("E24").Value = '  +1.97%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '0.164'
$ws.Range("E25").Value = '  +2.96%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '8.53'
$ws.Range("E26").Value = '  +5.50%  '

$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.56%  '

$ws.Range("E28").Value = '  +7.41%  '

$ws.Range("D29").Value = '0.0₃0820'
$ws.Range("E29").Value = '  +3.17%  '

$ws.Range("D30").Value = '6.84'
$ws.Range("E30").Value = '  +7.01%  '

$ws.Range("D31").Value = '168.88'
$ws.Range("E31").Value = '  +3.56%  '

$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("E33").Value = '  +3.04%  '

$ws.Range("D34").Value = '1.11'
$ws.Range("E34").Value = '  +11.85%  '

$ws.Range("E35").Value = '  +8.22%  '

$ws.Range("D36").Value = '1.37'
$ws.Range("E36").Value = '  +9.49%  '

$ws.Range("D37").Value = '1.71'
$ws.Range("E37").Value = '  +4.81%  '

$ws.Range("D38").Value = '338.45'
$ws.Range("E38").Value = '  +9.16%  '

$ws.Range("D39").Value = '0.940'
$ws.Range("E39").Value = '  +11.96%  '

$ws.Range("D40").Value = '4.13'
$ws.Range("E40").Value = '  +5.68%  '

$ws.Range("D41").Value = '38.35'
$ws.Range("E41").Value = '  +0.86%  '

$ws.Range("D42").Value = '5.35'
$ws.Range("E42").Value = '  +6.39%  '

$ws.Range("D43").Value = '0.0581'
$ws.Range("E43").Value = '  +5.86%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '21.19'
$ws.Range("E44").Value = '  +5.19%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.631'
$ws.Range("E45").Value = '  +4.53%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '135.79'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '20.57'
$ws.Range("E47").Value = '  +4.36%  '

$ws.Range("E48").Value = '  +4.65%  '

$ws.Range("E49").Value = '  +0.88%  '

$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("D51").Value = '2.091.39'
$ws.Range("E51").Value = '  +3.29%  '
